$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.043.44'
$ws.Range('E2').Value = '  +1.04%  '
$ws.Range('D3').Value = '2.297.26'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '507.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.14'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.995'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.531'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.25%  '
$ws.Range('D9').Value = '2.321.22'
$ws.Range('E9').Value = '  +0.85%  '
$ws.Range('E10').Value = '  +2.75%  '
$ws.Range('E11').Value = '  +1.71%  '
$ws.Range('E13').Value = '  +0.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.93'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = '2.707.28'
$ws.Range('E15').Value = '  +0.27%  '
$ws.Range('D16').Value = '55.062.65'
$ws.Range('E16').Value = '  +1.16%  '
$ws.Range('E17').Value = '  +1.62%  '
$ws.Range('D18').Value = '2.307.21'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.77'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.09%  '
$ws.Range('E20').Value = '  +0.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '310.98'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.64'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.08%  '
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.992'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '172.92'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('E29').Value = '  +3.09%  '
$ws.Range('D30').Value = '0.0₃0710'
$ws.Range('E30').Value = '  +2.24%  '
$ws.Range('E31').Value = '  +0.42%  '
$ws.Range('E32').Value = '  +5.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.10'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.993'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.23'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.921'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.83'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.84%  '
$ws.Range('E40').Value = '  +2.35%  '
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '135.17'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.01%  '
$ws.Range('E43').Value = '  +1.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.95'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '261.34'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.66%  '
$ws.Range('E46').Value = '  +2.11%  '
$ws.Range('E47').Value = '  +1.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.552'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.57%  '
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('E50').Value = '  +2.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.57'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.08%  '

Write-Host "Updated $($ws.Name) with latest crypto prices"
